# Updated cryptos list (Price + Volume(1h) columns, plus row 51 coin swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text so numeric-looking strings like
# "1.004" or "26.262.84" are not auto-converted to numbers by Excel,
# matching the inlineStr cells already used throughout this sheet.
function Set-CellText($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-CellText "D2" '26.262.84'
Set-CellText "D3" '1.657.37'
Set-CellText "E3" '  -0.97%  '
Set-CellText "D4" '1.004'
Set-CellText "D5" '219.56'
Set-CellText "E5" '  -0.81%  '
Set-CellText "D6" '0.5242'
Set-CellText "E6" '  -1.81%  '
Set-CellText "E7" '  -0.69%  '
Set-CellText "D8" '0.2676'
Set-CellText "E8" '  +0.27%  '
Set-CellText "D9" '0.06365'
Set-CellText "E9" '  -0.39%  '
Set-CellText "D10" '20.73'
Set-CellText "E10" '  -0.77%  '
Set-CellText "D11" '0.07706'
Set-CellText "E11" '  -1.94%  '
Set-CellText "D12" '4.608'
Set-CellText "E12" '  +1.60%  '
Set-CellText "D13" '1.617.90'
Set-CellText "E13" '  -2.95%  '
Set-CellText "D14" '1.886.18'
Set-CellText "E14" '  -0.89%  '
Set-CellText "D15" '0.5649'
Set-CellText "E15" '  +0.46%  '
Set-CellText "D16" '0.0₅8232'
Set-CellText "E16" '  +0.33%  '
Set-CellText "D17" '65.49'
Set-CellText "E17" '  -1.15%  '
Set-CellText "D18" '26.260.17'
Set-CellText "E18" '  -0.71%  '
Set-CellText "D19" '1.003'
Set-CellText "E19" '  -0.76%  '
Set-CellText "D20" '4.704'
Set-CellText "E20" '  -0.58%  '
Set-CellText "D21" '10.45'
Set-CellText "E21" '  +1.26%  '
Set-CellText "D22" '192.79'
Set-CellText "E22" '  -2.41%  '
Set-CellText "D23" '6.016'
Set-CellText "E23" '  -1.12%  '
Set-CellText "E24" '  -0.65%  '
Set-CellText "D25" '143.35'
Set-CellText "E25" '  -2.01%  '
Set-CellText "E26" '  -2.09%  '
Set-CellText "D27" '7.302'
Set-CellText "E27" '  +0.58%  '
Set-CellText "D28" '15.95'
Set-CellText "E28" '  -1.66%  '
Set-CellText "D29" '1.508'
Set-CellText "E29" '  -0.08%  '
Set-CellText "D30" '0.05655'
Set-CellText "E30" '  -4.62%  '
Set-CellText "E31" '  -1.13%  '
Set-CellText "D32" '3.513'
Set-CellText "E32" '  -1.47%  '
Set-CellText "D33" '3.357'
Set-CellText "E33" '  +0.48%  '
Set-CellText "E34" '  -1.67%  '
Set-CellText "E35" '  -1.25%  '
Set-CellText "D36" '0.9478'
Set-CellText "E36" '  -2.27%  '
Set-CellText "D37" '2.416'
Set-CellText "E37" '  -0.96%  '
Set-CellText "D38" '0.5777'
Set-CellText "E38" '  -1.15%  '
Set-CellText "D39" '0.01603'
Set-CellText "E39" '  -0.90%  '
Set-CellText "D40" '5.976'
Set-CellText "E40" '  +0.76%  '
Set-CellText "D41" '2.569'
Set-CellText "E41" '  -0.19%  '
Set-CellText "D42" '0.8467'
Set-CellText "E42" '  -2.23%  '
Set-CellText "E43" '  -0.72%  '
Set-CellText "D44" '1.021.85'
Set-CellText "E44" '  -5.46%  '
Set-CellText "D45" '101.68'
Set-CellText "E45" '  -1.57%  '
Set-CellText "D46" '1.796.78'
Set-CellText "D47" '58.51'
Set-CellText "E47" '  -0.30%  '
Set-CellText "E48" '  -0.15%  '
Set-CellText "D49" '1.004'
Set-CellText "E49" '  -0.97%  '
Set-CellText "D50" '0.05321'
Set-CellText "E50" '  +3.06%  '
Set-CellText "B51" 'Mantle'
Set-CellText "C51" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-CellText "D51" '0.4350'
Set-CellText "E51" '  -1.52%  '
